# The "Rules" worksheet cell E8 held the string "Good Morning" (shared
# string table entry) which is no longer referenced anywhere else in the
# sheet. Replacing it with the new text "GIT UPDATE" makes Excel drop the
# now-unused shared string and re-index the remaining ones, which is what
# produced the cascading <v> index shifts seen in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

# The diff also shows the sheet's active selection changed to E8.
$ws.Range("E8").Select()
